$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.011.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.33%  '
$ws.Range("D3").Value = "'1.962.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.27%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = "'326.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.37%  '
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").Value = "'0.4980"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.39%  '
$ws.Range("D8").Value = "'0.4196"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.47%  '
$ws.Range("D9").Value = "'52.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("D10").Value = "'0.09132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("E11").Value = '  -6.99%  '
$ws.Range("D12").Value = "'22.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.60%  '
$ws.Range("D13").Value = "'2.014.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = "'7.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.78%  '
$ws.Range("D15").Value = "'6.423"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.71%  '
$ws.Range("D16").Value = "'1.008"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.50%  '
$ws.Range("D18").Value = "'91.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -10.42%  '
$ws.Range("D19").Value = "'0.06665"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").Value = "'19.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.45%  '
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = "'5.959"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.13%  '
$ws.Range("D23").Value = "'29.045.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.21%  '
$ws.Range("D24").Value = "'12.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.94%  '
$ws.Range("D25").Value = "'2.283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").Value = "'2.226.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("D27").Value = "'156.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.97%  '
$ws.Range("D28").Value = "'20.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.90%  '
$ws.Range("D29").Value = "'6.145"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.48%  '
$ws.Range("D30").Value = "'2.254"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.77%  '
$ws.Range("D31").Value = "'126.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.30%  '
$ws.Range("D32").Value = "'1.035"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.93%  '
$ws.Range("D33").Value = "'0.09818"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.55%  '
$ws.Range("D34").Value = "'1.520"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.19%  '
$ws.Range("D35").Value = "'5.747"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.37%  '
$ws.Range("D36").Value = "'3.676"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.88%  '
$ws.Range("D37").Value = "'0.02406"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.31%  '
$ws.Range("D38").Value = "'1.301"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("D39").Value = "'0.06310"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.69%  '
$ws.Range("D40").Value = "'8.896"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.06%  '
$ws.Range("D41").Value = "'0.6419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.96%  '
$ws.Range("D42").Value = "'11.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.89%  '
$ws.Range("D43").Value = "'0.1975"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.54%  '
$ws.Range("D44").Value = "'1.006"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").Value = "'0.6202"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.41%  '
$ws.Range("D46").Value = "'13.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.46%  '
$ws.Range("D47").Value = "'2.169"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.05%  '
$ws.Range("D48").Value = "'1.289"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = "'3.461"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.64%  '
$ws.Range("D50").Value = "'0.00000000331"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").Value = "'0.06942"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.45%  '
